# Add a new worksheet named "L6" at the end of the workbook (after the
# last existing sheet) containing the "last six games" summary table:
# Form / Goals scored / Goals conceded / Total Goals per team.

$wb = $excel.ActiveWorkbook

$origActive = $wb.ActiveSheet
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)

$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "L6"

# Header row
$ws.Range("B1").Value = "Form"
$ws.Range("C1").Value = "Goals scored"
$ws.Range("D1").Value = "Goals conceded"
$ws.Range("E1").Value = "Total Goals"

# Keep column A ("1".."10") as text, matching the rest of the workbook
# (e.g. the "Table" sheet's rank column), not auto-converted to numbers.
$ws.Range("A2:A11").NumberFormat = "@"

$data = @(
    ,@("1", "Alloa,D L L L D W", "Alloa,2 1 0 1 1 1", "Alloa,2 2 6 2 1 0", "Alloa,4 3 6 3 2 1")
    ,@("2", "Arbroath,L W D W L D", "Arbroath,0 2 2 4 3 0", "Arbroath,1 1 2 0 4 0", "Arbroath,1 3 4 4 7 0")
    ,@("3", "Ayr,D L D L D D", "Ayr,2 0 1 0 0 2", "Ayr,2 3 1 4 0 2", "Ayr,4 3 2 4 0 4")
    ,@("4", "Dundee,W D D D W W", "Dundee,3 1 0 1 2 2", "Dundee,0 1 0 1 1 0", "Dundee,3 2 0 2 3 2")
    ,@("5", "Dunfermline,D D D W W L", "Dunfermline,0 1 0 3 4 0", "Dunfermline,0 1 0 1 3 1", "Dunfermline,0 2 0 4 7 1")
    ,@("6", "Hearts,L D W D W W", "Hearts,2 0 6 0 3 4", "Hearts,3 0 0 0 0 0", "Hearts,5 0 6 0 3 4")
    ,@("7", "Inverness C,W W D D L D", "Inverness C,1 4 1 1 0 2", "Inverness C,0 1 1 1 3 2", "Inverness C,1 5 2 2 3 4")
    ,@("8", "Morton,L L D D D D", "Morton,0 1 1 0 1 0", "Morton,1 4 1 0 1 0", "Morton,1 5 2 0 2 0")
    ,@("9", "Queen of Sth,W D L L D L", "Queen of Sth,3 1 0 1 0 0", "Queen of Sth,2 1 1 3 0 2", "Queen of Sth,5 2 1 4 0 2")
    ,@("10", "Raith Rvs,W D W W L L", "Raith Rvs,5 2 1 2 1 0", "Raith Rvs,1 2 0 1 2 4", "Raith Rvs,6 4 1 3 3 4")
)

$row = 2
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $ws.Cells.Item($row, 5).Value = $r[4]
    $row = $row + 1
}

$origActive.Activate()
